# Applies the "Add files via upload" edit: append 30 new guild-war
# result rows (rows 64-93) to Sheet1, matching the author's new data,
# and update the sheet view's selection to reflect where they were
# editing (N86) plus scroll the visible window down to around row 46.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(64, 1).Value = "플라튼, 트루드, 아멜리아"
$ws.Cells.Item(64, 2).Value = "루"
$ws.Cells.Item(64, 3).Value = "아1아2트2"
$ws.Cells.Item(64, 4).Value = "크리스, 로지, 녹스"
$ws.Cells.Item(64, 5).Value = "맬패로"
$ws.Cells.Item(64, 6).Value = "크1크2녹2"
$ws.Cells.Item(64, 7).Value = "선"
$ws.Cells.Item(64, 8).Value = 260105
$ws.Cells.Item(64, 9).Value = "모현"
$ws.Cells.Item(64, 10).Value = "공격"
$ws.Cells.Item(65, 1).Value = "프레이야, 바네사, 연희"
$ws.Cells.Item(65, 2).Value = "연지"
$ws.Cells.Item(65, 3).Value = "바1프2키2"
$ws.Cells.Item(65, 4).Value = "프레이야, 바네사, 키리엘"
$ws.Cells.Item(65, 5).Value = "연지"
$ws.Cells.Item(65, 6).Value = "프2바1연2"
$ws.Cells.Item(65, 7).Value = "후"
$ws.Cells.Item(65, 8).Value = 260105
$ws.Cells.Item(65, 9).Value = "모현"
$ws.Cells.Item(65, 10).Value = "공격"
$ws.Cells.Item(66, 1).Value = "카일, 카구라, 파이"
$ws.Cells.Item(66, 2).Value = "이린"
$ws.Cells.Item(66, 3).Value = "구2카1카2"
$ws.Cells.Item(66, 4).Value = "겔리두스, 엘리시아, 트루드"
$ws.Cells.Item(66, 5).Value = "루"
$ws.Cells.Item(66, 6).Value = "겔2트2엘1"
$ws.Cells.Item(66, 7).Value = "후"
$ws.Cells.Item(66, 8).Value = 260106
$ws.Cells.Item(66, 9).Value = "모현"
$ws.Cells.Item(66, 10).Value = "공격"
$ws.Cells.Item(67, 1).Value = "카일, 카구라, 델론즈"
$ws.Cells.Item(67, 2).Value = "이린"
$ws.Cells.Item(67, 3).Value = "구2카1카2"
$ws.Cells.Item(67, 4).Value = "프레이야, 바네사, 밀리아"
$ws.Cells.Item(67, 5).Value = "연지"
$ws.Cells.Item(67, 6).Value = "바1밀2프2"
$ws.Cells.Item(67, 7).Value = "선"
$ws.Cells.Item(67, 8).Value = 260107
$ws.Cells.Item(67, 9).Value = "모현"
$ws.Cells.Item(67, 10).Value = "공격"
$ws.Cells.Item(68, 1).Value = "플라튼, 콜트, 엘리시아"
$ws.Cells.Item(68, 2).Value = "루"
$ws.Cells.Item(68, 3).Value = "엘2엘1콜2"
$ws.Cells.Item(68, 4).Value = "플라튼, 아멜리아, 트루드"
$ws.Cells.Item(68, 5).Value = "루"
$ws.Cells.Item(68, 6).Value = "아2트2트1"
$ws.Cells.Item(68, 7).Value = "후"
$ws.Cells.Item(68, 8).Value = 260108
$ws.Cells.Item(68, 9).Value = "모현"
$ws.Cells.Item(68, 10).Value = "공격"
$ws.Cells.Item(69, 1).Value = "카일, 카구라, 파이"
$ws.Cells.Item(69, 2).Value = "이린"
$ws.Cells.Item(69, 3).Value = "파1카1카2"
$ws.Cells.Item(69, 4).Value = "프레이야, 밀리아, 바네사"
$ws.Cells.Item(69, 5).Value = "연지"
$ws.Cells.Item(69, 6).Value = "프2바1프1"
$ws.Cells.Item(69, 7).Value = "후"
$ws.Cells.Item(69, 8).Value = 260109
$ws.Cells.Item(69, 9).Value = "모현"
$ws.Cells.Item(69, 10).Value = "공격"
$ws.Cells.Item(70, 1).Value = "카일, 카구라, 파이"
$ws.Cells.Item(70, 2).Value = "이린"
$ws.Cells.Item(70, 3).Value = "카1파1카2"
$ws.Cells.Item(70, 4).Value = "트루드, 겔리두스, 엘리시아"
$ws.Cells.Item(70, 5).Value = "루"
$ws.Cells.Item(70, 6).Value = "엘1겔2트2"
$ws.Cells.Item(70, 7).Value = "후"
$ws.Cells.Item(70, 8).Value = 260110
$ws.Cells.Item(70, 9).Value = "모현"
$ws.Cells.Item(70, 10).Value = "공격"
$ws.Cells.Item(71, 1).Value = "카일, 카구라, 파이"
$ws.Cells.Item(71, 2).Value = "이린"
$ws.Cells.Item(71, 3).Value = "구2카1카2"
$ws.Cells.Item(71, 4).Value = "카일, 카구라, 파이"
$ws.Cells.Item(71, 5).Value = "이린"
$ws.Cells.Item(71, 6).Value = "구2카1카2"
$ws.Cells.Item(71, 7).Value = "후"
$ws.Cells.Item(71, 8).Value = 260111
$ws.Cells.Item(71, 9).Value = "모현"
$ws.Cells.Item(71, 10).Value = "공격"
$ws.Cells.Item(72, 1).Value = "프레이야, 밀리아, 바네사"
$ws.Cells.Item(72, 2).Value = "연지"
$ws.Cells.Item(72, 3).Value = "밀2프2바1"
$ws.Cells.Item(72, 4).Value = "오공, 겔리두스, 스파이크"
$ws.Cells.Item(72, 5).Value = "파이크"
$ws.Cells.Item(72, 6).Value = "오2겔2스2"
$ws.Cells.Item(72, 7).Value = "후"
$ws.Cells.Item(72, 8).Value = 260112
$ws.Cells.Item(72, 9).Value = "모현"
$ws.Cells.Item(72, 10).Value = "공격"
$ws.Cells.Item(73, 1).Value = "플라튼, 프레이야, 콜트"
$ws.Cells.Item(73, 2).Value = "카람"
$ws.Cells.Item(73, 3).Value = "프2콜1콜2"
$ws.Cells.Item(73, 4).Value = "카일, 파이, 아멜리아"
$ws.Cells.Item(73, 5).Value = "이린"
$ws.Cells.Item(73, 6).Value = "아2카1카2"
$ws.Cells.Item(73, 7).Value = "선"
$ws.Cells.Item(73, 8).Value = 260113
$ws.Cells.Item(73, 9).Value = "모현"
$ws.Cells.Item(73, 10).Value = "공격"
$ws.Cells.Item(74, 1).Value = "플라튼, 트루드, 아멜리아"
$ws.Cells.Item(74, 2).Value = "루"
$ws.Cells.Item(74, 3).Value = "아1아2트2"
$ws.Cells.Item(74, 4).Value = "카일, 카구라, 파이"
$ws.Cells.Item(74, 5).Value = "이린"
$ws.Cells.Item(74, 6).Value = "구1카1카2"
$ws.Cells.Item(74, 7).Value = "선"
$ws.Cells.Item(74, 8).Value = 260114
$ws.Cells.Item(74, 9).Value = "모현"
$ws.Cells.Item(74, 10).Value = "공격"
$ws.Cells.Item(75, 1).Value = "밀리아, 멜키르, 프레이야"
$ws.Cells.Item(75, 2).Value = "유"
$ws.Cells.Item(75, 3).Value = "밀2프2멜2"
$ws.Cells.Item(75, 4).Value = "프레이야, 쥬리, 바네사"
$ws.Cells.Item(75, 5).Value = "연지"
$ws.Cells.Item(75, 6).Value = "바1쥬2프2"
$ws.Cells.Item(75, 7).Value = "후"
$ws.Cells.Item(75, 8).Value = 260115
$ws.Cells.Item(75, 9).Value = "모현"
$ws.Cells.Item(75, 10).Value = "공격"
$ws.Cells.Item(76, 1).Value = "프레이야, 멜키르, 밀리아"
$ws.Cells.Item(76, 2).Value = "연지"
$ws.Cells.Item(76, 3).Value = "멜2프2밀2"
$ws.Cells.Item(76, 4).Value = "프레이야, 엘리시아, 밀리아"
$ws.Cells.Item(76, 5).Value = "연지"
$ws.Cells.Item(76, 6).Value = "엘1프2프1"
$ws.Cells.Item(76, 7).Value = "후"
$ws.Cells.Item(76, 8).Value = 260116
$ws.Cells.Item(76, 9).Value = "모현"
$ws.Cells.Item(76, 10).Value = "공격"
$ws.Cells.Item(77, 1).Value = "플라튼, 콜트, 프레이야"
$ws.Cells.Item(77, 2).Value = "루"
$ws.Cells.Item(77, 3).Value = "프2콜1콜2"
$ws.Cells.Item(77, 4).Value = "트루드, 오공, 엘리시아"
$ws.Cells.Item(77, 5).Value = "크리"
$ws.Cells.Item(77, 6).Value = "오2엘1오1"
$ws.Cells.Item(77, 7).Value = "선"
$ws.Cells.Item(77, 8).Value = 260117
$ws.Cells.Item(77, 9).Value = "모현"
$ws.Cells.Item(77, 10).Value = "공격"
$ws.Cells.Item(78, 1).Value = "밀리아, 연희, 멜키르"
$ws.Cells.Item(78, 2).Value = "크리"
$ws.Cells.Item(78, 3).Value = "멜1멜2연2"
$ws.Cells.Item(78, 4).Value = "트루드, 스파이크, 오공"
$ws.Cells.Item(78, 5).Value = "루"
$ws.Cells.Item(78, 6).Value = "오2트2스2"
$ws.Cells.Item(78, 7).Value = "후"
$ws.Cells.Item(78, 8).Value = 260118
$ws.Cells.Item(78, 9).Value = "모현"
$ws.Cells.Item(78, 10).Value = "공격"
$ws.Cells.Item(79, 1).Value = "플라튼, 실베스타, 아멜리아"
$ws.Cells.Item(79, 2).Value = "루"
$ws.Cells.Item(79, 3).Value = "아1실2실1"
$ws.Cells.Item(79, 4).Value = "크리스, 녹스, 로지"
$ws.Cells.Item(79, 5).Value = "맬패로"
$ws.Cells.Item(79, 6).Value = "크1로1녹2"
$ws.Cells.Item(79, 7).Value = "선"
$ws.Cells.Item(79, 8).Value = 260119
$ws.Cells.Item(79, 9).Value = "모현"
$ws.Cells.Item(79, 10).Value = "공격"
$ws.Cells.Item(80, 1).Value = "에이스, 콜트, 파이"
$ws.Cells.Item(80, 2).Value = "리첼"
$ws.Cells.Item(80, 3).Value = "파1콜1콜2"
$ws.Cells.Item(80, 4).Value = "카일, 카구라, 파이"
$ws.Cells.Item(80, 5).Value = "이린"
$ws.Cells.Item(80, 6).Value = "구2카1카2"
$ws.Cells.Item(80, 7).Value = "선"
$ws.Cells.Item(80, 8).Value = 260120
$ws.Cells.Item(80, 9).Value = "모현"
$ws.Cells.Item(80, 10).Value = "공격"
$ws.Cells.Item(81, 1).Value = "프레이야, 바네사, 키리엘"
$ws.Cells.Item(81, 2).Value = "연지"
$ws.Cells.Item(81, 3).Value = "바1프2프1"
$ws.Cells.Item(81, 4).Value = "오공, 겔리두스, 엘리시아"
$ws.Cells.Item(81, 5).Value = "루"
$ws.Cells.Item(81, 6).Value = "겔2오2겔1"
$ws.Cells.Item(81, 7).Value = "후"
$ws.Cells.Item(81, 8).Value = 260121
$ws.Cells.Item(81, 9).Value = "모현"
$ws.Cells.Item(81, 10).Value = "공격"
$ws.Cells.Item(82, 1).Value = "프레이야, 바네사, 키리엘"
$ws.Cells.Item(82, 2).Value = "연지"
$ws.Cells.Item(82, 3).Value = "바1프2프1"
$ws.Cells.Item(82, 4).Value = "오공, 겔리두스, 스파이크"
$ws.Cells.Item(82, 5).Value = "파이크"
$ws.Cells.Item(82, 6).Value = "오2겔2겔1"
$ws.Cells.Item(82, 7).Value = "후"
$ws.Cells.Item(82, 8).Value = 260122
$ws.Cells.Item(82, 9).Value = "모현"
$ws.Cells.Item(82, 10).Value = "공격"
$ws.Cells.Item(83, 1).Value = "프레이야, 멜키르, 밀리아"
$ws.Cells.Item(83, 2).Value = "연지"
$ws.Cells.Item(83, 3).Value = "밀2멜2프2"
$ws.Cells.Item(83, 4).Value = "프레이야, 멜키르, 밀리아"
$ws.Cells.Item(83, 5).Value = "리첼"
$ws.Cells.Item(83, 6).Value = "밀2프2프1"
$ws.Cells.Item(83, 7).Value = "후"
$ws.Cells.Item(83, 8).Value = 260123
$ws.Cells.Item(83, 9).Value = "모현"
$ws.Cells.Item(83, 10).Value = "공격"
$ws.Cells.Item(84, 1).Value = "프레이야, 바네사, 키리엘"
$ws.Cells.Item(84, 2).Value = "연지"
$ws.Cells.Item(84, 3).Value = "바1프2키1"
$ws.Cells.Item(84, 4).Value = "프레이야, 밀리아, 바네사"
$ws.Cells.Item(84, 5).Value = "연지"
$ws.Cells.Item(84, 6).Value = "바1밀2프2"
$ws.Cells.Item(84, 7).Value = "선"
$ws.Cells.Item(84, 8).Value = 260124
$ws.Cells.Item(84, 9).Value = "모현"
$ws.Cells.Item(84, 10).Value = "공격"
$ws.Cells.Item(85, 1).Value = "오공, 스파이크, 로지"
$ws.Cells.Item(85, 2).Value = "파이크"
$ws.Cells.Item(85, 3).Value = "오2오1스2"
$ws.Cells.Item(85, 4).Value = "오공, 스파이크, 플라튼"
$ws.Cells.Item(85, 5).Value = "카람"
$ws.Cells.Item(85, 6).Value = "오2플2오1"
$ws.Cells.Item(85, 7).Value = "선"
$ws.Cells.Item(85, 8).Value = 260125
$ws.Cells.Item(85, 9).Value = "모현"
$ws.Cells.Item(85, 10).Value = "공격"
$ws.Cells.Item(86, 1).Value = "연희, 바네사, 키리엘"
$ws.Cells.Item(86, 2).Value = "연지"
$ws.Cells.Item(86, 3).Value = "바1키2연2"
$ws.Cells.Item(86, 4).Value = "카일, 카구라, 파이"
$ws.Cells.Item(86, 5).Value = "이린"
$ws.Cells.Item(86, 6).Value = "구2파2카2"
$ws.Cells.Item(86, 7).Value = "선"
$ws.Cells.Item(86, 8).Value = 260126
$ws.Cells.Item(86, 9).Value = "모현"
$ws.Cells.Item(86, 10).Value = "공격"
$ws.Cells.Item(87, 1).Value = "프레이야, 바네사, 키리엘"
$ws.Cells.Item(87, 2).Value = "연지"
$ws.Cells.Item(87, 3).Value = "바1프2프1"
$ws.Cells.Item(87, 4).Value = "플라튼, 트루드, 아멜리아"
$ws.Cells.Item(87, 5).Value = "루"
$ws.Cells.Item(87, 6).Value = "아2아1트2"
$ws.Cells.Item(87, 7).Value = "후"
$ws.Cells.Item(87, 8).Value = 260127
$ws.Cells.Item(87, 9).Value = "모현"
$ws.Cells.Item(87, 10).Value = "공격"
$ws.Cells.Item(88, 1).Value = "프레이야, 바네사, 키리엘"
$ws.Cells.Item(88, 2).Value = "연지"
$ws.Cells.Item(88, 3).Value = "바1프2프1"
$ws.Cells.Item(88, 4).Value = "오공, 겔리두스, 엘리시아"
$ws.Cells.Item(88, 5).Value = "카람"
$ws.Cells.Item(88, 6).Value = "오2겔2겔1"
$ws.Cells.Item(88, 7).Value = "후"
$ws.Cells.Item(88, 8).Value = 260128
$ws.Cells.Item(88, 9).Value = "모현"
$ws.Cells.Item(88, 10).Value = "공격"
$ws.Cells.Item(89, 1).Value = "연희, 바네사, 키리엘"
$ws.Cells.Item(89, 2).Value = "연지"
$ws.Cells.Item(89, 3).Value = "바1연2키2"
$ws.Cells.Item(89, 4).Value = "연희, 니아, 로지"
$ws.Cells.Item(89, 5).Value = "맬패로"
$ws.Cells.Item(89, 6).Value = "니2연2니1"
$ws.Cells.Item(89, 7).Value = "선"
$ws.Cells.Item(89, 8).Value = 260129
$ws.Cells.Item(89, 9).Value = "모현"
$ws.Cells.Item(89, 10).Value = "공격"
$ws.Cells.Item(90, 1).Value = "카일, 카구라, 파이"
$ws.Cells.Item(90, 2).Value = "이린"
$ws.Cells.Item(90, 3).Value = "파1카1카2"
$ws.Cells.Item(90, 4).Value = "프레이야, 바네사, 쥬리"
$ws.Cells.Item(90, 5).Value = "연지"
$ws.Cells.Item(90, 6).Value = "바1쥬2프1"
$ws.Cells.Item(90, 7).Value = "선"
$ws.Cells.Item(90, 8).Value = 260130
$ws.Cells.Item(90, 9).Value = "모현"
$ws.Cells.Item(90, 10).Value = "공격"
$ws.Cells.Item(91, 1).Value = "프레이야, 멜키르, 밀리아"
$ws.Cells.Item(91, 2).Value = "이린"
$ws.Cells.Item(91, 3).Value = "멜2프2밀2"
$ws.Cells.Item(91, 4).Value = "카일, 엘리시아, 카구라"
$ws.Cells.Item(91, 5).Value = "이린"
$ws.Cells.Item(91, 6).Value = "엘1카1구2"
$ws.Cells.Item(91, 7).Value = "선"
$ws.Cells.Item(91, 8).Value = 260131
$ws.Cells.Item(91, 9).Value = "모현"
$ws.Cells.Item(91, 10).Value = "공격"
$ws.Cells.Item(92, 1).Value = "밀리아, 멜키르, 프레이야"
$ws.Cells.Item(92, 2).Value = "연지"
$ws.Cells.Item(92, 3).Value = "멜2프2밀2"
$ws.Cells.Item(92, 4).Value = "밀리아, 바네사, 프레이야"
$ws.Cells.Item(92, 5).Value = "루"
$ws.Cells.Item(92, 6).Value = "바1프2프1"
$ws.Cells.Item(92, 7).Value = "선"
$ws.Cells.Item(92, 8).Value = 260132
$ws.Cells.Item(92, 9).Value = "모현"
$ws.Cells.Item(92, 10).Value = "공격"
$ws.Cells.Item(93, 1).Value = "트루드, 겔리두스, 엘리시아"
$ws.Cells.Item(93, 2).Value = "유"
$ws.Cells.Item(93, 3).Value = "트2엘1겔2"
$ws.Cells.Item(93, 4).Value = "스파이크, 엘리스, 리나"
$ws.Cells.Item(93, 5).Value = "파이크"
$ws.Cells.Item(93, 6).Value = "엘1스2"
$ws.Cells.Item(93, 7).Value = "후"
$ws.Cells.Item(93, 8).Value = 260133
$ws.Cells.Item(93, 9).Value = "모현"
$ws.Cells.Item(93, 10).Value = "공격"

# Match the style (center-aligned, like the rest of the table) used by
# the pre-existing data rows.
$newRange = $ws.Range("A64:J93")
$newRange.HorizontalAlignment = -4108
$newRange.VerticalAlignment = -4108

# Reflect the author's final cursor position / scroll offset captured
# in the saved sheetView.
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("N86").Select()
